$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Queue Report")

# Clear the old dummy report rows (2:3) entirely before laying down the
# new "graphical representation" data (3 rows of queue-item summaries).
$ws.Range("A2:I3").ClearContents()

$itemDate = 43870.092361111114
$subDate  = 43869.863194444442
$dateFmt  = "m/d/yy h:mm"

for ($r = 2; $r -le 4; $r++) {
    $ws.Cells.Item($r, 1).Value = $itemDate
    $ws.Cells.Item($r, 1).NumberFormat = $dateFmt

    $ws.Cells.Item($r, 2).Value = "Successful"

    $ws.Cells.Item($r, 3).Value = $subDate
    $ws.Cells.Item($r, 3).NumberFormat = $dateFmt

    $ws.Cells.Item($r, 5).Value = $subDate
    $ws.Cells.Item($r, 5).NumberFormat = $dateFmt

    $ws.Cells.Item($r, 9).Value = "ED"
}
